# Updates Leve-profit market-data columns (H-N) across sheets, matching
# scheduled-runner market refresh. Plain values only (no formulas in source).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 206
$ws.Range("J9").Value = 300
$ws.Range("L9").Value = 300
$ws.Range("N9").Value = -638
$ws.Range("H32").Value = 4421.8
$ws.Range("I32").Value = 4658.1665
$ws.Range("K32").Value = 4658.1665
$ws.Range("M32").Value = -4332.1665
$ws.Range("H92").Value = 1059.8
$ws.Range("I92").Value = 448.8889
$ws.Range("J92").Value = 1976.1666
$ws.Range("K92").Value = 448.8889
$ws.Range("L92").Value = 1976.1666
$ws.Range("M92").Value = 799.1111000000001
$ws.Range("N92").Value = -4472.1666
$ws.Range("H132").Value = 1452.2
$ws.Range("I132").Value = 1239.6086
$ws.Range("K132").Value = 3718.8258
$ws.Range("M132").Value = -1188.8258
$ws.Range("H138").Value = 2153.6323
$ws.Range("J138").Value = 2297.087
$ws.Range("L138").Value = 6891.261
$ws.Range("N138").Value = -17171.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10096.667
$ws.Range("I61").Value = 9578.235000000001
$ws.Range("J61").Value = 12300
$ws.Range("K61").Value = 9578.235000000001
$ws.Range("L61").Value = 12300
$ws.Range("M61").Value = -9366.235000000001
$ws.Range("N61").Value = -12724
$ws.Range("H74").Value = 4345.0586
$ws.Range("I74").Value = 3181.0908
$ws.Range("J74").Value = 6479
$ws.Range("K74").Value = 3181.0908
$ws.Range("L74").Value = 6479
$ws.Range("M74").Value = -2307.0908
$ws.Range("N74").Value = -8227
$ws.Range("H77").Value = 4345.0586
$ws.Range("I77").Value = 3181.0908
$ws.Range("J77").Value = 6479
$ws.Range("K77").Value = 15905.454
$ws.Range("L77").Value = 32395
$ws.Range("M77").Value = -11537.454
$ws.Range("N77").Value = -41131
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H132").Value = 7201.0967
$ws.Range("I132").Value = 6149.923
$ws.Range("K132").Value = 18449.769
$ws.Range("M132").Value = -15919.769
$ws.Range("H136").Value = 10096.667
$ws.Range("I136").Value = 9578.235000000001
$ws.Range("J136").Value = 12300
$ws.Range("K136").Value = 28734.705
$ws.Range("L136").Value = 36900
$ws.Range("M136").Value = -26184.705
$ws.Range("N136").Value = -42000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3186.4688
$ws.Range("I20").Value = 2321.4736
$ws.Range("J20").Value = 4450.6924
$ws.Range("K20").Value = 2321.4736
$ws.Range("L20").Value = 4450.6924
$ws.Range("M20").Value = -2074.4736
$ws.Range("N20").Value = -4944.6924
$ws.Range("H86").Value = 1333.1904
$ws.Range("I86").Value = 1333.1904
$ws.Range("K86").Value = 1333.1904
$ws.Range("M86").Value = -210.1904
$ws.Range("H89").Value = 1333.1904
$ws.Range("I89").Value = 1333.1904
$ws.Range("K89").Value = 6665.951999999999
$ws.Range("M89").Value = -1049.951999999999
$ws.Range("H134").Value = 3885.4849
$ws.Range("I134").Value = 3885.4849
$ws.Range("K134").Value = 11656.4547
$ws.Range("M134").Value = -9121.4547

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6929
$ws.Range("J99").Value = 6781.6665
$ws.Range("L99").Value = 6781.6665
$ws.Range("N99").Value = -9777.666499999999
$ws.Range("H107").Value = 520.94446
$ws.Range("I107").Value = 352.2857
$ws.Range("K107").Value = 352.2857
$ws.Range("M107").Value = 1567.7143
$ws.Range("H126").Value = 6929
$ws.Range("J126").Value = 6781.6665
$ws.Range("L126").Value = 20344.9995
$ws.Range("N126").Value = -25284.9995
$ws.Range("H132").Value = 3493.625
$ws.Range("I132").Value = 3391.9333
$ws.Range("J132").Value = 3663.111
$ws.Range("K132").Value = 10175.7999
$ws.Range("L132").Value = 10989.333
$ws.Range("M132").Value = -7645.7999
$ws.Range("N132").Value = -16049.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 388499.5
$ws.Range("I21").Value = 511999.34
$ws.Range("K21").Value = 511999.34
$ws.Range("M21").Value = -511826.34
$ws.Range("H30").Value = 388499.5
$ws.Range("I30").Value = 511999.34
$ws.Range("K30").Value = 511999.34
$ws.Range("M30").Value = -511894.34
$ws.Range("H39").Value = 47750
$ws.Range("I39").Value = 46000
$ws.Range("K39").Value = 46000
$ws.Range("M39").Value = -45468
$ws.Range("H126").Value = 2505897.2
$ws.Range("I126").Value = 4004336.8
$ws.Range("K126").Value = 12013010.4
$ws.Range("M126").Value = -12010540.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3169.9
$ws.Range("I7").Value = 2957.1428
$ws.Range("K7").Value = 2957.1428
$ws.Range("M7").Value = -2845.1428
$ws.Range("H16").Value = 1102.409
$ws.Range("I16").Value = 941.7222
$ws.Range("J16").Value = 1825.5
$ws.Range("K16").Value = 941.7222
$ws.Range("L16").Value = 1825.5
$ws.Range("M16").Value = -771.7222
$ws.Range("N16").Value = -2165.5
$ws.Range("H46").Value = 12740.25
$ws.Range("I46").Value = 5354.727
$ws.Range("K46").Value = 5354.727
$ws.Range("M46").Value = -5166.727
$ws.Range("H122").Value = 2909.125
$ws.Range("I122").Value = 2596.4666
$ws.Range("J122").Value = 3430.2222
$ws.Range("K122").Value = 7789.399800000001
$ws.Range("L122").Value = 10290.6666
$ws.Range("M122").Value = -5339.399800000001
$ws.Range("N122").Value = -15190.6666
$ws.Range("H126").Value = 3169.9
$ws.Range("I126").Value = 2957.1428
$ws.Range("K126").Value = 8871.428400000001
$ws.Range("M126").Value = -6401.428400000001
$ws.Range("H136").Value = 4562.423
$ws.Range("I136").Value = 4330
$ws.Range("J136").Value = 5538.6
$ws.Range("K136").Value = 12990
$ws.Range("L136").Value = 16615.8
$ws.Range("M136").Value = -10440
$ws.Range("N136").Value = -21715.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 19762.834
$ws.Range("I74").Value = 15979
$ws.Range("J74").Value = 20519.6
$ws.Range("K74").Value = 15979
$ws.Range("L74").Value = 20519.6
$ws.Range("M74").Value = -15043
$ws.Range("N74").Value = -22391.6
$ws.Range("H77").Value = 19762.834
$ws.Range("I77").Value = 15979
$ws.Range("J77").Value = 20519.6
$ws.Range("K77").Value = 47937
$ws.Range("L77").Value = 61558.8
$ws.Range("M77").Value = -43257
$ws.Range("N77").Value = -70918.79999999999
$ws.Range("H81").Value = 2043.4286
$ws.Range("J81").Value = 4891
$ws.Range("L81").Value = 9782
$ws.Range("N81").Value = -11904
$ws.Range("H84").Value = 2043.4286
$ws.Range("J84").Value = 4891
$ws.Range("L84").Value = 48910
$ws.Range("N84").Value = -59518
$ws.Range("H96").Value = 2637.1333
$ws.Range("I96").Value = 2228.6667
$ws.Range("K96").Value = 2228.6667
$ws.Range("M96").Value = -855.6667000000002
$ws.Range("H122").Value = 4071.6206
$ws.Range("I122").Value = 4040.5
$ws.Range("J122").Value = 4140.778
$ws.Range("K122").Value = 12121.5
$ws.Range("L122").Value = 12422.334
$ws.Range("M122").Value = -9671.5
$ws.Range("N122").Value = -17322.334

